$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price (D) and Volume(1h) (E) columns hold plain text values in the
# source data (several "prices" even contain multiple "." thousands
# separators, e.g. "23.524.06", so they are not valid numbers at all).
# Force the target range to Text format before writing so Excel's
# automatic type inference does not silently convert values such as
# "51.20" into the number 51.2. Afterwards restore the cell style to
# "Normal" so no stray formatting is introduced (we only needed the
# number format long enough to perform the assignment).
$priceVolRange = $ws.Range("D2:E51")
$priceVolRange.NumberFormat = "@"

$ws.Range("D2").Value = "23.524.06"
$ws.Range("E2").Value = "  +1.15%  "
$ws.Range("D3").Value = "1.653.59"
$ws.Range("E3").Value = "  +2.44%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("E5").Value = "  +0.02%  "
$ws.Range("D6").Value = "302.29"
$ws.Range("E6").Value = "  -0.19%  "
$ws.Range("D7").Value = "0.3834"
$ws.Range("E7").Value = "  +1.47%  "
$ws.Range("D8").Value = "51.20"
$ws.Range("E8").Value = "  -0.94%  "
$ws.Range("D9").Value = "0.3592"
$ws.Range("E9").Value = "  +1.80%  "
$ws.Range("D10").Value = "0.08196"
$ws.Range("E10").Value = "  +1.21%  "
$ws.Range("D11").Value = "1.239"
$ws.Range("E11").Value = "  +2.68%  "
$ws.Range("D12").Value = "1.001"
$ws.Range("E12").Value = "  -0.04%  "
$ws.Range("E13").Value = "  +0.92%  "
$ws.Range("D14").Value = "6.474"
$ws.Range("E14").Value = "  +1.69%  "
$ws.Range("D15").Value = "7.489"
$ws.Range("E15").Value = "  +3.04%  "
$ws.Range("E16").Value = "  +0.80%  "
$ws.Range("D17").Value = "1.648.63"
$ws.Range("E17").Value = "  +3.94%  "
$ws.Range("D18").Value = "97.44"
$ws.Range("E18").Value = "  +3.75%  "
$ws.Range("D19").Value = "0.06975"
$ws.Range("E19").Value = "  +1.13%  "
$ws.Range("D20").Value = "6.813"
$ws.Range("E20").Value = "  +5.16%  "
$ws.Range("D21").Value = "17.65"
$ws.Range("E21").Value = "  +2.53%  "
$ws.Range("D22").Value = "1.002"
$ws.Range("E22").Value = "  +0.03%  "
$ws.Range("D23").Value = "12.67"
$ws.Range("E23").Value = "  +2.83%  "
$ws.Range("D24").Value = "23.527.84"
$ws.Range("E24").Value = "  +1.26%  "
$ws.Range("D25").Value = "2.511"
$ws.Range("E25").Value = "  -0.07%  "
$ws.Range("D26").Value = "3.014"
$ws.Range("E26").Value = "  +0.09%  "
$ws.Range("D27").Value = "21.19"
$ws.Range("E27").Value = "  +1.54%  "
$ws.Range("D28").Value = "152.10"
$ws.Range("E28").Value = "  +0.66%  "
$ws.Range("D29").Value = "5.242"
$ws.Range("E29").Value = "  +0.00%  "
$ws.Range("D30").Value = "133.70"
$ws.Range("E30").Value = "  +1.24%  "
$ws.Range("B31").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C31").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D31").Value = "1.836.67"
$ws.Range("E31").Value = "  +3.60%  "
$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").Value = "7.203"
$ws.Range("E32").Value = "  +11.50%  "
$ws.Range("D33").Value = "2.253"
$ws.Range("E33").Value = "  +7.70%  "
$ws.Range("E34").Value = "  +6.92%  "
$ws.Range("D35").Value = "1.058"
$ws.Range("E35").Value = "  -0.69%  "
$ws.Range("D36").Value = "0.02804"
$ws.Range("E36").Value = "  +3.70%  "
$ws.Range("D37").Value = "6.117"
$ws.Range("E37").Value = "  +4.63%  "
$ws.Range("D38").Value = "0.2497"
$ws.Range("E38").Value = "  +1.80%  "
$ws.Range("D39").Value = "0.08782"
$ws.Range("E39").Value = "  +1.15%  "
$ws.Range("D40").Value = "0.07016"
$ws.Range("E40").Value = "  +1.35%  "
$ws.Range("D41").Value = "13.22"
$ws.Range("E41").Value = "  +11.03%  "
$ws.Range("D42").Value = "0.7007"
$ws.Range("E42").Value = "  +1.85%  "
$ws.Range("D43").Value = "1.335"
$ws.Range("E43").Value = "  +0.70%  "
$ws.Range("E44").Value = "  +5.28%  "
$ws.Range("D45").Value = "0.6517"
$ws.Range("E45").Value = "  +3.17%  "
$ws.Range("E46").Value = "  +0.09%  "
$ws.Range("D47").Value = "2.305"
$ws.Range("E47").Value = "  +2.44%  "
$ws.Range("D48").Value = "3.955"
$ws.Range("E48").Value = "  +0.34%  "
$ws.Range("D49").Value = "0.07904"
$ws.Range("E49").Value = "  +0.48%  "
$ws.Range("D50").Value = "128.03"
$ws.Range("E50").Value = "  +0.34%  "
$ws.Range("D51").Value = "1.192"
$ws.Range("E51").Value = "  +1.91%  "

$priceVolRange.Style = "Normal"
